# Update TPM-derived NATMI metrics for the Adam9-Itga6 ligand-receptor sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 11.92133566666667
$ws.Range("H2").Value = 35.764007
$ws.Range("I2").Value = 0.10300114145944
$ws.Range("J2").Value = 0.10300114145944
$ws.Range("M2").Value = 173.5452066666667
$ws.Range("N2").Value = 520.63562
$ws.Range("O2").Value = 0.6098887991422922
$ws.Range("P2").Value = 0.6098887991422922
$ws.Range("Q2").Value = 2068.890662014371
$ws.Range("R2").Value = 18620.01595812934
$ws.Range("S2").Value = 0.06281924247498324
$ws.Range("T2").Value = 0.06281924247498323

# Row 3
$ws.Range("G3").Value = 11.92133566666667
$ws.Range("H3").Value = 35.764007
$ws.Range("I3").Value = 0.10300114145944
$ws.Range("J3").Value = 0.10300114145944
$ws.Range("O3").Value = 0.003264284357140855
$ws.Range("P3").Value = 0.003264284357140855
$ws.Range("Q3").Value = 11.07324389978289
$ws.Range("R3").Value = 99.659195098046
$ws.Range("S3").Value = 0.0003362250148337024
$ws.Range("T3").Value = 0.0003362250148337024

# Row 4
$ws.Range("G4").Value = 11.92133566666667
$ws.Range("H4").Value = 35.764007
$ws.Range("I4").Value = 0.10300114145944
$ws.Range("J4").Value = 0.10300114145944
$ws.Range("M4").Value = 54.64271666666667
$ws.Range("N4").Value = 163.92815
$ws.Range("O4").Value = 0.192030546333187
$ws.Range("P4").Value = 0.192030546333187
$ws.Range("Q4").Value = 651.4141671218945
$ws.Range("R4").Value = 5862.727504097051
$ws.Range("S4").Value = 0.01977936546739815
$ws.Range("T4").Value = 0.01977936546739814

# Row 5
$ws.Range("G5").Value = 11.92133566666667
$ws.Range("H5").Value = 35.764007
$ws.Range("I5").Value = 0.10300114145944
$ws.Range("J5").Value = 0.10300114145944
$ws.Range("M5").Value = 1.069012
$ws.Range("N5").Value = 3.207036
$ws.Range("O5").Value = 0.00375682196858928
$ws.Range("P5").Value = 0.00375682196858928
$ws.Range("Q5").Value = 12.74405088369467
$ws.Range("R5").Value = 114.696457953252
$ws.Range("S5").Value = 0.0003869569510245963
$ws.Range("T5").Value = 0.0003869569510245964

# Row 6
$ws.Range("G6").Value = 11.92133566666667
$ws.Range("H6").Value = 35.764007
$ws.Range("I6").Value = 0.10300114145944
$ws.Range("J6").Value = 0.10300114145944
$ws.Range("M6").Value = 54.36641700000001
$ws.Range("N6").Value = 163.099251
$ws.Range("O6").Value = 0.1910595481987908
$ws.Range("P6").Value = 0.1910595481987908
$ws.Range("Q6").Value = 648.120306050973
$ws.Range("R6").Value = 5833.082754458757
$ws.Range("S6").Value = 0.01967935155120034
$ws.Range("T6").Value = 0.01967935155120034

# Row 7
$ws.Range("G7").Value = 42.86866899999999
$ws.Range("I7").Value = 0.3703881816022666
$ws.Range("J7").Value = 0.3703881816022665
$ws.Range("M7").Value = 173.5452066666667
$ws.Range("N7").Value = 520.63562
$ws.Range("O7").Value = 0.6098887991422922
$ws.Range("P7").Value = 0.6098887991422922
$ws.Range("Q7").Value = 7439.652021129925
$ws.Range("R7").Value = 66956.86819016933
$ws.Range("S7").Value = 0.2258956032939036
$ws.Range("T7").Value = 0.2258956032939036

# Row 8
$ws.Range("G8").Value = 42.86866899999999
$ws.Range("I8").Value = 0.3703881816022666
$ws.Range("J8").Value = 0.3703881816022665
$ws.Range("O8").Value = 0.003264284357140855
$ws.Range("P8").Value = 0.003264284357140855
$ws.Range("Q8").Value = 39.81896330822732
$ws.Range("S8").Value = 0.001209052347274125
$ws.Range("T8").Value = 0.001209052347274125

# Row 9
$ws.Range("G9").Value = 42.86866899999999
$ws.Range("I9").Value = 0.3703881816022666
$ws.Range("J9").Value = 0.3703881816022665
$ws.Range("M9").Value = 54.64271666666667
$ws.Range("N9").Value = 163.92815
$ws.Range("O9").Value = 0.192030546333187
$ws.Range("P9").Value = 0.192030546333187
$ws.Range("Q9").Value = 2342.460534044116
$ws.Range("R9").Value = 21082.14480639705
$ws.Range("S9").Value = 0.07112584486843894
$ws.Range("T9").Value = 0.07112584486843893

# Row 10
$ws.Range("G10").Value = 42.86866899999999
$ws.Range("I10").Value = 0.3703881816022666
$ws.Range("J10").Value = 0.3703881816022665
$ws.Range("M10").Value = 1.069012
$ws.Range("N10").Value = 3.207036
$ws.Range("O10").Value = 0.00375682196858928
$ws.Range("P10").Value = 0.00375682196858928
$ws.Range("Q10").Value = 45.827121585028
$ws.Range("R10").Value = 412.444094265252
$ws.Range("S10").Value = 0.001391482457549231
$ws.Range("T10").Value = 0.001391482457549231

# Row 11
$ws.Range("G11").Value = 42.86866899999999
$ws.Range("I11").Value = 0.3703881816022666
$ws.Range("J11").Value = 0.3703881816022665
$ws.Range("M11").Value = 54.36641700000001
$ws.Range("N11").Value = 163.099251
$ws.Range("O11").Value = 0.1910595481987908
$ws.Range("P11").Value = 0.1910595481987908
$ws.Range("Q11").Value = 2330.615935088972
$ws.Range("R11").Value = 20975.54341580076
$ws.Range("S11").Value = 0.07076619863510071
$ws.Range("T11").Value = 0.0707661986351007

# Row 12
$ws.Range("G12").Value = 27.63817166666666
$ws.Range("H12").Value = 82.91451499999999
$ws.Range("I12").Value = 0.23879566091562
$ws.Range("J12").Value = 0.23879566091562
$ws.Range("M12").Value = 173.5452066666667
$ws.Range("N12").Value = 520.63562
$ws.Range("O12").Value = 0.6098887991422922
$ws.Range("P12").Value = 0.6098887991422922
$ws.Range("Q12").Value = 4796.472213780477
$ws.Range("R12").Value = 43168.2499240243
$ws.Range("S12").Value = 0.1456387988762175
$ws.Range("T12").Value = 0.1456387988762175

# Row 13
$ws.Range("G13").Value = 27.63817166666666
$ws.Range("H13").Value = 82.91451499999999
$ws.Range("I13").Value = 0.23879566091562
$ws.Range("J13").Value = 0.23879566091562
$ws.Range("O13").Value = 0.003264284357140855
$ws.Range("P13").Value = 0.003264284357140855
$ws.Range("Q13").Value = 25.67197370885222
$ws.Range("R13").Value = 231.04776337967
$ws.Range("S13").Value = 0.0007794969404799703
$ws.Range("T13").Value = 0.0007794969404799701

# Row 14
$ws.Range("G14").Value = 27.63817166666666
$ws.Range("H14").Value = 82.91451499999999
$ws.Range("I14").Value = 0.23879566091562
$ws.Range("J14").Value = 0.23879566091562
$ws.Range("M14").Value = 54.64271666666667
$ws.Range("N14").Value = 163.92815
$ws.Range("O14").Value = 0.192030546333187
$ws.Range("P14").Value = 0.192030546333187
$ws.Range("Q14").Value = 1510.224783566361
$ws.Range("R14").Value = 13592.02305209725
$ws.Range("S14").Value = 0.045856061227621
$ws.Range("T14").Value = 0.04585606122762098

# Row 15
$ws.Range("G15").Value = 27.63817166666666
$ws.Range("H15").Value = 82.91451499999999
$ws.Range("I15").Value = 0.23879566091562
$ws.Range("J15").Value = 0.23879566091562
$ws.Range("M15").Value = 1.069012
$ws.Range("N15").Value = 3.207036
$ws.Range("O15").Value = 0.00375682196858928
$ws.Range("P15").Value = 0.00375682196858928
$ws.Range("Q15").Value = 29.54553716972667
$ws.Range("R15").Value = 265.90983452754
$ws.Range("S15").Value = 0.0008971127849315977
$ws.Range("T15").Value = 0.0008971127849315977

# Row 16
$ws.Range("G16").Value = 27.63817166666666
$ws.Range("H16").Value = 82.91451499999999
$ws.Range("I16").Value = 0.23879566091562
$ws.Range("J16").Value = 0.23879566091562
$ws.Range("M16").Value = 54.36641700000001
$ws.Range("N16").Value = 163.099251
$ws.Range("O16").Value = 0.1910595481987908
$ws.Range("P16").Value = 0.1910595481987908
$ws.Range("Q16").Value = 1502.588365947585
$ws.Range("R16").Value = 13523.29529352826
$ws.Range("S16").Value = 0.04562419108637
$ws.Range("T16").Value = 0.04562419108636999

# Row 17
$ws.Range("G17").Value = 6.093442
$ws.Range("H17").Value = 18.280326
$ws.Range("I17").Value = 0.05264774845421206
$ws.Range("J17").Value = 0.05264774845421205
$ws.Range("M17").Value = 173.5452066666667
$ws.Range("N17").Value = 520.63562
$ws.Range("O17").Value = 0.6098887991422922
$ws.Range("P17").Value = 0.6098887991422922
$ws.Range("Q17").Value = 1057.487651201347
$ws.Range("R17").Value = 9517.38886081212
$ws.Range("S17").Value = 0.03210927208228487
$ws.Range("T17").Value = 0.03210927208228486

# Row 18
$ws.Range("G18").Value = 6.093442
$ws.Range("H18").Value = 18.280326
$ws.Range("I18").Value = 0.05264774845421206
$ws.Range("J18").Value = 0.05264774845421205
$ws.Range("O18").Value = 0.003264284357140855
$ws.Range("P18").Value = 0.003264284357140855
$ws.Range("Q18").Value = 5.659950473825333
$ws.Range("R18").Value = 50.939554264428
$ws.Range("S18").Value = 0.0001718572217177711
$ws.Range("T18").Value = 0.000171857221717771

# Row 19
$ws.Range("G19").Value = 6.093442
$ws.Range("H19").Value = 18.280326
$ws.Range("I19").Value = 0.05264774845421206
$ws.Range("J19").Value = 0.05264774845421205
$ws.Range("M19").Value = 54.64271666666667
$ws.Range("N19").Value = 163.92815
$ws.Range("O19").Value = 0.192030546333187
$ws.Range("P19").Value = 0.192030546333187
$ws.Range("Q19").Value = 332.9622247307666
$ws.Range("R19").Value = 2996.6600225769
$ws.Range("S19").Value = 0.01010997589887455
$ws.Range("T19").Value = 0.01010997589887454

# Row 20
$ws.Range("G20").Value = 6.093442
$ws.Range("H20").Value = 18.280326
$ws.Range("I20").Value = 0.05264774845421206
$ws.Range("J20").Value = 0.05264774845421205
$ws.Range("M20").Value = 1.069012
$ws.Range("N20").Value = 3.207036
$ws.Range("O20").Value = 0.00375682196858928
$ws.Range("P20").Value = 0.00375682196858928
$ws.Range("Q20").Value = 6.513962619304
$ws.Range("R20").Value = 58.62566357373601
$ws.Range("S20").Value = 0.0001977882179895462
$ws.Range("T20").Value = 0.0001977882179895462

# Row 21
$ws.Range("G21").Value = 6.093442
$ws.Range("H21").Value = 18.280326
$ws.Range("I21").Value = 0.05264774845421206
$ws.Range("J21").Value = 0.05264774845421205
$ws.Range("M21").Value = 54.36641700000001
$ws.Range("N21").Value = 163.099251
$ws.Range("O21").Value = 0.1910595481987908
$ws.Range("P21").Value = 0.1910595481987908
$ws.Range("Q21").Value = 331.278608737314
$ws.Range("R21").Value = 2981.507478635826
$ws.Range("S21").Value = 0.01005885503334534
$ws.Range("T21").Value = 0.01005885503334534

# Row 22
$ws.Range("G22").Value = 27.218222
$ws.Range("H22").Value = 81.65466599999999
$ws.Range("I22").Value = 0.2351672675684614
$ws.Range("J22").Value = 0.2351672675684614
$ws.Range("M22").Value = 173.5452066666667
$ws.Range("N22").Value = 520.63562
$ws.Range("O22").Value = 0.6098887991422922
$ws.Range("P22").Value = 0.6098887991422922
$ws.Range("Q22").Value = 4723.591962089213
$ws.Range("R22").Value = 42512.32765880292
$ws.Range("S22").Value = 0.1434258824149031
$ws.Range("T22").Value = 0.143425882414903

# Row 23
$ws.Range("G23").Value = 27.218222
$ws.Range("H23").Value = 81.65466599999999
$ws.Range("I23").Value = 0.2351672675684614
$ws.Range("J23").Value = 0.2351672675684614
$ws.Range("O23").Value = 0.003264284357140855
$ws.Range("P23").Value = 0.003264284357140855
$ws.Range("Q23").Value = 25.28189954143867
$ws.Range("R23").Value = 227.537095872948
$ws.Range("S23").Value = 0.0007676528328352865
$ws.Range("T23").Value = 0.0007676528328352864

# Row 24
$ws.Range("G24").Value = 27.218222
$ws.Range("H24").Value = 81.65466599999999
$ws.Range("I24").Value = 0.2351672675684614
$ws.Range("J24").Value = 0.2351672675684614
$ws.Range("M24").Value = 54.64271666666667
$ws.Range("N24").Value = 163.92815
$ws.Range("O24").Value = 0.192030546333187
$ws.Range("P24").Value = 0.192030546333187
$ws.Range("Q24").Value = 1487.277592916433
$ws.Range("R24").Value = 13385.4983362479
$ws.Range("S24").Value = 0.04515929887085442
$ws.Range("T24").Value = 0.04515929887085442

# Row 25
$ws.Range("G25").Value = 27.218222
$ws.Range("H25").Value = 81.65466599999999
$ws.Range("I25").Value = 0.2351672675684614
$ws.Range("J25").Value = 0.2351672675684614
$ws.Range("M25").Value = 1.069012
$ws.Range("N25").Value = 3.207036
$ws.Range("O25").Value = 0.00375682196858928
$ws.Range("P25").Value = 0.00375682196858928
$ws.Range("Q25").Value = 29.096605936664
$ws.Range("R25").Value = 261.869453429976
$ws.Range("S25").Value = 0.0008834815570943091
$ws.Range("T25").Value = 0.0008834815570943091

# Row 26
$ws.Range("G26").Value = 27.218222
$ws.Range("H26").Value = 81.65466599999999
$ws.Range("I26").Value = 0.2351672675684614
$ws.Range("J26").Value = 0.2351672675684614
$ws.Range("M26").Value = 54.36641700000001
$ws.Range("N26").Value = 163.099251
$ws.Range("O26").Value = 0.1910595481987908
$ws.Range("P26").Value = 0.1910595481987908
$ws.Range("Q26").Value = 1479.757207250574
$ws.Range("R26").Value = 13317.81486525517
$ws.Range("S26").Value = 0.04493095189277438
$ws.Range("T26").Value = 0.04493095189277437

